$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "June" -> "July"
$ws.Range("H4").Value = "July"

# Clear a handful of stray promo/discount figures
$ws.Range("D13").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("J38").Value = ""
$ws.Range("D42").Value = ""

# New model "Z30" inserted into the sorted list (row 41), pushing
# Z25 and Z50 down by one row each.
$ws.Range("G41").Value = "Z30"
$ws.Range("H41").Value = 9300
$ws.Range("I41").Value = 9790

$ws.Range("G42").Value = "Z25"
$ws.Range("H42").Value = 8310
$ws.Range("I42").Value = 8990

$ws.Range("G43").Value = "Z50"
$ws.Range("H43").Value = 10340
$ws.Range("I43").Value = 10990
$ws.Range("K43").Value = 0

# Update selection / scroll position
$ws.Range("P13").Select()
